$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price and 1h volume change values for rows 2-44
# D-column (Price) values are prefixed with a leading apostrophe so Excel
# stores them as literal text (matching the original text-formatted cells)
# instead of auto-converting look-alike numbers to the Number type.
$ws.Range("D2").Value = "'30.112.30"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "'1.909.97"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'320.31"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.5062"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.4082"
$ws.Range("E8").Value = "  +4.56%  "
$ws.Range("D9").Value = "'0.08346"
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("D10").Value = "'42.40"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "'1.105"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "'23.78"
$ws.Range("E12").Value = "  +5.47%  "
$ws.Range("D13").Value = "'1.909.71"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "'6.398"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "'7.234"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'92.42"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "'0.00001097"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("D19").Value = "'0.06497"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").Value = "'18.46"
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'5.939"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").Value = "'30.127.37"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").Value = "'2.127.66"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "'21.78"
$ws.Range("E27").Value = "  +4.84%  "
$ws.Range("D28").Value = "'162.94"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'2.283"
$ws.Range("E29").Value = "  +3.44%  "
$ws.Range("D30").Value = "'128.72"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").Value = "'1.143"
$ws.Range("E31").Value = "  +10.26%  "
$ws.Range("D32").Value = "'0.1043"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "'5.960"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("D34").Value = "'3.786"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").Value = "'0.02456"
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("D36").Value = "'5.353"
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("D37").Value = "'0.06394"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "'0.2149"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").Value = "'0.6524"
$ws.Range("E39").Value = "  +4.45%  "
$ws.Range("D40").Value = "'1.196"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").Value = "'8.614"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").Value = "'11.41"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").Value = "'1.213"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "'13.38"
$ws.Range("E44").Value = "  +4.52%  "

# Row 45 and 46 swap: NEARProtocol/Decentraland order reversed, with updated price/volume data
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6082"
$ws.Range("E45").Value = "  +3.82%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.193"
$ws.Range("E46").Value = "  +10.68%  "

# Rows 47-51 price/volume updates
$ws.Range("D47").Value = "'3.621"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "'1.209"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").Value = "'122.05"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'79.06"
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("D51").Value = "'1.137"
$ws.Range("E51").Value = "  +1.58%  "
